$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new hours-log entry (row 20) for Erica.
# Column A holds the date as literal text (matching the existing rows, which
# are plain strings rather than real dates), so build it via a text formula
# and then convert it to a plain value in place -- this avoids Excel's
# automatic "looks like a date" conversion (which would turn it into a date
# serial number + a new date-formatted style) while still leaving the cell
# as an ordinary, unstyled text value like the rows above it.
$ws.Range("A20").Formula = '="2/6/2010"'
$ws.Range("A20").Copy()
$ws.Range("A20").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("B20").Value = 3
$ws.Range("C20").Value = "CFP Update and misc"

# Move the active selection down to the next empty row, as in the source file
$ws.Range("A21").Select()
